# train.xlsx refresh: re-captured traffic sample (celery fix / load_app func).
# The recorded API-traffic log gained a re-run at 2025-03-06 08:54:05, several
# request indices got renumbered, a few endpoint/flag values were corrected,
# and the old trailing row (13) is dropped from the sheet entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row outright (was row 11 / index 13). Deleting the row
# (rather than just clearing it) shifts nothing else and shrinks the sheet's
# dimension from A1:Q11 down to A1:Q10, matching the target layout.
$ws.Rows.Item(11).Delete()

# Row 2
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "2025-03-06 08:54:05"

# Row 3
$ws.Range("D3").Value = "2025-03-06 08:54:05"

# Row 4
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = "2025-03-06 08:54:05"
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = $false

# Row 5
$ws.Range("A5").Value = 9
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = "2025-03-06 08:54:05"
$ws.Range("F5").Value = "http://49.234.6.241:5230/api/v1/memo/21"
$ws.Range("G5").Value = "/api/v1/memo/21"

# Row 6
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = "2025-03-06 08:54:05"
$ws.Range("F6").Value = "http://47.97.114.24:5230/api/v1/memo/21"
$ws.Range("G6").Value = "/api/v1/memo/21"

# Row 7
$ws.Range("A7").Value = 11
$ws.Range("C7").Value = 11
$ws.Range("D7").Value = "2025-03-06 08:54:05"

# Row 8
$ws.Range("A8").Value = 12
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = "2025-03-06 08:54:05"
$ws.Range("F8").Value = "http://47.97.114.24:5230/api/v1/resource/16"
$ws.Range("G8").Value = "/api/v1/resource/16"
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $true
$ws.Range("Q8").Value = $true

# Row 9
$ws.Range("A9").Value = 13
$ws.Range("C9").Value = 13
$ws.Range("D9").Value = "2025-03-06 08:54:05"
$ws.Range("M9").Value = 0.003
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = $false

# Row 10
$ws.Range("A10").Value = 14
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 14
$ws.Range("D10").Value = "2025-03-06 08:54:05"
$ws.Range("M10").Value = 0.004
